$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Plan1")

# Correção de funções: substitui o valor "?" por "!" nas células de status
# das linhas correspondentes.
$ws.Range("F22").Value = "!"
$ws.Range("F23").Value = "!"
$ws.Range("F24").Value = "!"
$ws.Range("F25").Value = "!"
$ws.Range("F29").Value = "!"
$ws.Range("F35").Value = "!"
$ws.Range("F39").Value = "!"

# Ajusta a seleção/posição da visualização para a célula F46
$ws.Range("F46").Select()
